$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.977.44"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.820.49"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.08"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4678"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3662"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07351"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8739"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.29"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "1.803.45"
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.418"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07148"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.59"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008747"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "27.006.99"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.292"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "2.052.61"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.889"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.04"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.40"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.141"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.239"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.78"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08881"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7560"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.159"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.501"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.945"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05311"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01946"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.973"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.383"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.180"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5297"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1652"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.462"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4891"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.46"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.665"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.10"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06294"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.13%  "
